# Turns the single-sheet workbook into a two-sheet ValidLogin / InvalidLogin
# login-data workbook:
#   - Sheet1 is renamed to "ValidLogin" and filled with a valid UserName/
#     Password pair (admin / manager).
#   - A new sheet "InvalidLogin" is added right after it, filled with an
#     invalid UserName/Password pair (Bhanu / Damager), column widths sized
#     to fit the data, and left as the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- ValidLogin (renamed original sheet) ---------------------------------
$validLogin = $wb.Worksheets.Item(1)
$validLogin.Name = "ValidLogin"

$validLogin.Range("A1").Value = "UserName"
$validLogin.Range("B1").Value = "Password"
$validLogin.Range("A2").Value = "admin"
$validLogin.Range("B2").Value = "manager"

# --- InvalidLogin (new sheet, placed after ValidLogin) --------------------
$invalidLogin = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $validLogin)
$invalidLogin.Name = "InvalidLogin"

$invalidLogin.Range("A1").Value = "UserName"
$invalidLogin.Range("B1").Value = "Password"
$invalidLogin.Range("A2").Value = "Bhanu"
$invalidLogin.Range("B2").Value = "Damager"

# Size the columns to fit their (longer) content, like the real workbook's
# best-fit column widths.
$invalidLogin.Columns.Item(1).ColumnWidth = 10.28515625
$invalidLogin.Columns.Item(2).ColumnWidth = 9.42578125

$invalidLogin.Range("B3").Select() | Out-Null

# --- Views: selection + zoom on each sheet, InvalidLogin ends up active ---
$validLogin.Range("A1:B2").Select() | Out-Null
$validLogin.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 235

$invalidLogin.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 250
